$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.817.22"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.318.84"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "2.676.45"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "2.323.63"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "42.747.17"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +14.45%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "147.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0700"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +22.10%  "
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").Value = "1.921.07"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.75%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").Value = "2.546.89"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.60%  "
